$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: median_rent_price / Richmond, VA ---
$ws.Range("A26").Value = 22
$ws.Range("B26").Value = "median_rent_price"
$ws.Range("C26").Value = "Richmond, VA median rent price for year"
$ws.Range("D26").Value = "int"
$ws.Range("E26").Value = "Quant"
$ws.Range("F26").Value = "Discrete"
$ws.Range("G26").Value = "Dollars"
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 2006
$ws.Range("M26").Value = 2017
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = 1
$ws.Range("S26").Value = 0
$ws.Range("T26").Value = 0
$ws.Range("U26").Value = 0
$ws.Range("V26").Value = "location_name == Richmond, VA"
$ws.Range("W26").Value = 944
$ws.Range("X26").Value = 1060
$ws.Range("Y26").Value = "https://www.deptofnumbers.com/rent/virginia/richmond/"
$ws.Range("Z26").Value = "web scraped to CSV"
$ws.Range("AA26").Value = "../DataSet/"
$ws.Range("AB26").Value = "Dept of Numbers"
$ws.Range("AC26").Value = "median gross rent"
$ws.Range("AD26").Value = "No"

# --- Row 27: mean_rent_price / Richmond, VA ---
$ws.Range("A27").Value = 23
$ws.Range("B27").Value = "mean_rent_price"
$ws.Range("C27").Value = "Richmond, VA mean rent price for year"
$ws.Range("D27").Value = "int"
$ws.Range("E27").Value = "Quant"
$ws.Range("F27").Value = "Discrete"
$ws.Range("G27").Value = "Dollars"
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2006
$ws.Range("M27").Value = 2017
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = 1
$ws.Range("S27").Value = 0
$ws.Range("T27").Value = 0
$ws.Range("U27").Value = 0
$ws.Range("V27").Value = "location_name == Richmond, VA"
$ws.Range("W27").Value = 922
$ws.Range("X27").Value = 1059
$ws.Range("Y27").Value = "https://www.deptofnumbers.com/rent/virginia/richmond/"
$ws.Range("Z27").Value = "web scraped to CSV"
$ws.Range("AA27").Value = "../DataSet/"
$ws.Range("AB27").Value = "Dept of Numbers"
$ws.Range("AC27").Value = "mean gross rent"
$ws.Range("AD27").Value = "No"

# --- Row 28: duplicate of row 26, to be re-pointed at Virginia Beach ---
$ws.Range("A28").Value = 22
$ws.Range("B28").Value = "median_rent_price"
$ws.Range("C28").Value = "Richmond, VA median rent price for year"
$ws.Range("D28").Value = "int"
$ws.Range("E28").Value = "Quant"
$ws.Range("F28").Value = "Discrete"
$ws.Range("G28").Value = "Dollars"
$ws.Range("H28").Value = 1
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 2006
$ws.Range("M28").Value = 2017
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = 1
$ws.Range("S28").Value = 0
$ws.Range("T28").Value = 0
$ws.Range("U28").Value = 0
$ws.Range("V28").Value = "location_name == Richmond, VA"
$ws.Range("W28").Value = 944
$ws.Range("X28").Value = 1060
$ws.Range("Y28").Value = "https://www.deptofnumbers.com/rent/virginia/richmond/"
$ws.Range("Z28").Value = "web scraped to CSV"
$ws.Range("AA28").Value = "../DataSet/"
$ws.Range("AB28").Value = "Dept of Numbers"
$ws.Range("AC28").Value = "median gross rent"
$ws.Range("AD28").Value = "No"

# --- Row 29: duplicate of row 27, to be re-pointed at Virginia Beach ---
$ws.Range("A29").Value = 23
$ws.Range("B29").Value = "mean_rent_price"
$ws.Range("C29").Value = "Richmond, VA mean rent price for year"
$ws.Range("D29").Value = "int"
$ws.Range("E29").Value = "Quant"
$ws.Range("F29").Value = "Discrete"
$ws.Range("G29").Value = "Dollars"
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 2006
$ws.Range("M29").Value = 2017
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("R29").Value = 1
$ws.Range("S29").Value = 0
$ws.Range("T29").Value = 0
$ws.Range("U29").Value = 0
$ws.Range("V29").Value = "location_name == Richmond, VA"
$ws.Range("W29").Value = 922
$ws.Range("X29").Value = 1059
$ws.Range("Y29").Value = "https://www.deptofnumbers.com/rent/virginia/richmond/"
$ws.Range("Z29").Value = "web scraped to CSV"
$ws.Range("AA29").Value = "../DataSet/"
$ws.Range("AB29").Value = "Dept of Numbers"
$ws.Range("AC29").Value = "mean gross rent"
$ws.Range("AD29").Value = "No"

# Fix up the descriptions for the Virginia Beach rows ...
$ws.Range("C28").Value = "Virginia Beach, VA median rent price for year"
$ws.Range("C29").Value = "Virginia Beach, VA mean rent price for year"

# ... and then the spatial identifier filter for both Virginia Beach rows.
$ws.Range("V28").Value = "location_name == Virginia Beach, VA"
$ws.Range("V29").Value = "location_name == Virginia Beach, VA"

# Leave the active selection on the last-edited cell, matching the author's
# final cursor position before saving.
[void]$ws.Range("W28").Select()
